# "Vary rain by month" — the Climate sheet's rain.rate column (D) no longer
# holds a flat 0.09 for every month; give each month its own observed rate,
# and tidy up the display so the differing precision reads clearly.

$wb = $excel.ActiveWorkbook
$climate = $wb.Worksheets.Item("Climate")

# New monthly rain.rate values (D2:D6), replacing the constant 0.09.
$climate.Range("D2").Value = 0.074
$climate.Range("D3").Value = 0.06
$climate.Range("D4").Value = 0.066
$climate.Range("D5").Value = 0.1
$climate.Range("D6").Value = 0.13

# D3 now reads to 3 decimal places, D5 to 2 — the other rows keep the
# sheet's default (general) formatting.
$climate.Range("D3").NumberFormat = "0.000"
$climate.Range("D5").NumberFormat = "0.00"

# The author ended the edit on the Climate tab (previously Application was
# active), having last clicked on E10.
$climate.Activate()
$climate.Range("E10").Select()
